$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet: Summary
# ---------------------------------------------------------------------
$wsSummary = $wb.Worksheets.Item("Summary")
$wsSummary.Range("B3").Value = 1400.54   # Current Capital
$wsSummary.Range("B4").Value = 0.34      # Total P&L $
$wsSummary.Range("B5").Value = 0.09      # Total P&L %
$wsSummary.Range("B6").Value = 76        # Total Trades
$wsSummary.Range("B7").Value = 35        # Winning Trades
$wsSummary.Range("B9").Value = 46.05     # Win Rate %

# ---------------------------------------------------------------------
# Sheet: Strategy Status (MarketMaking row = row 5)
# ---------------------------------------------------------------------
$wsStatus = $wb.Worksheets.Item("Strategy Status")
$wsStatus.Range("C5").Value = 100.54
$wsStatus.Range("D5").Value = 43
$wsStatus.Range("E5").Value = 0.23
$wsStatus.Range("F5").Value = 0.54
$wsStatus.Range("G5").Value = 48.84

# ---------------------------------------------------------------------
# Sheet: All Trades - update trade #48 / row 77 (now CLOSED)
# ---------------------------------------------------------------------
$wsAll = $wb.Worksheets.Item("All Trades")
$wsAll.Range("G77").Value = 0.17
$wsAll.Range("H77").Value = "CLOSED"
$wsAll.Range("I77").Value = 21.4286
$wsAll.Range("J77").Value = 0.03
$wsAll.Range("K77").Value = 100.54
$wsAll.Range("L77").Value = "early_exit"
$wsAll.Range("M77").Value = 0.14

# Append new trade #109 as row 110
$wsAll.Range("A110").Value = 109
$wsAll.Range("B110").NumberFormat = "@"
$wsAll.Range("B110").Value = "2026-02-17"
$wsAll.Range("B110").ClearFormats()
$wsAll.Range("C110").Value = "21:07:41"
$wsAll.Range("D110").Value = "MarketMaking"
$wsAll.Range("E110").Value = "DOWN"
$wsAll.Range("F110").Value = 0.14
$wsAll.Range("H110").Value = "OPEN"
$wsAll.Range("I110").Value = 0
$wsAll.Range("J110").Value = 0
$wsAll.Range("K110").Value = 100.5114872031006
$wsAll.Range("M110").Value = 0
$wsAll.Range("N110").Value = 0
$wsAll.Range("O110").Value = 0
$wsAll.Range("P110").Value = 0.6
$wsAll.Range("Q110").Value = "Normal spread capture: 19600 bps"

# ---------------------------------------------------------------------
# Sheet: MarketMaking - update trade #48 / row 44 (now CLOSED)
# ---------------------------------------------------------------------
$wsMM = $wb.Worksheets.Item("MarketMaking")
$wsMM.Range("G44").Value = 0.17
$wsMM.Range("H44").Value = "CLOSED"
$wsMM.Range("I44").Value = 21.4286
$wsMM.Range("J44").Value = 0.03
$wsMM.Range("K44").Value = 100.54
$wsMM.Range("P44").Value = "early_exit"
$wsMM.Range("Q44").Value = 0.14

# Append new trade #109 as row 77
$wsMM.Range("A77").Value = 109
$wsMM.Range("B77").NumberFormat = "@"
$wsMM.Range("B77").Value = "2026-02-17"
$wsMM.Range("B77").ClearFormats()
$wsMM.Range("C77").Value = "21:07:41"
$wsMM.Range("D77").Value = "MarketMaking"
$wsMM.Range("E77").Value = "DOWN"
$wsMM.Range("F77").Value = 0.14
$wsMM.Range("H77").Value = "OPEN"
$wsMM.Range("I77").Value = 0
$wsMM.Range("J77").Value = 0
$wsMM.Range("K77").Value = 100.5114872031006
$wsMM.Range("L77").Value = 0
$wsMM.Range("M77").Value = 0
$wsMM.Range("N77").Value = 0.6
$wsMM.Range("O77").Value = "Normal spread capture: 19600 bps"
$wsMM.Range("Q77").Value = 0

$wb.Save()

Write-Host "All edits applied successfully"
